$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status for the 54172188... row (row 3) changes from "Ready for handoff" to
# "Handback transform failed" on every sheet that shows it.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New "Error Detail" (column L) entries for the 54172188... row (row 3),
# reporting the handback/handoff file name mismatch per locale.
$wsZhCn.Range("L3").Value = "Handback file name: 00sckfi0.rht is different with handoff file name: 54172188-5415-447a-9f09-5a6684877614.be3761a2ce0a166bacf6d18629bbbd359b7a320c.zh-cn."
$wsDeDe.Range("L3").Value = "Handback file name: 00sckfi0.rht is different with handoff file name: 54172188-5415-447a-9f09-5a6684877614.be3761a2ce0a166bacf6d18629bbbd359b7a320c.de-de."
